# This script updates the "想去人数" (F column) figures for a handful of
# events on both the "展览" sheet and the "全部类型" sheet, matching the
# refreshed data snapshot described by the commit
# "Update gh-pages to output generated at 456a3b4".

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (row -> new F value) ---
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F6").Value  = 14202
$wsExpo.Range("F7").Value  = 15986
$wsExpo.Range("F9").Value  = 69
$wsExpo.Range("F24").Value = 6359
$wsExpo.Range("F26").Value = 1107
$wsExpo.Range("F27").Value = 5637
$wsExpo.Range("F31").Value = 4637

# --- Sheet "全部类型" (row -> new F value) ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F6").Value  = 14202
$wsAll.Range("F7").Value  = 15986
$wsAll.Range("F9").Value  = 69
$wsAll.Range("F25").Value = 6359
$wsAll.Range("F27").Value = 1107
$wsAll.Range("F29").Value = 5637
$wsAll.Range("F33").Value = 4637
